$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "Time" column (old V),
# shifting old V,W,X,Y -> X,Y,Z,AA
$ws.Columns("V:W").Insert()

# New header cells
$ws.Range("V1").Value = "Second"
$ws.Range("W1").Value = "P"

# New "Second" column: convert FinishTime (U) fraction-of-day into whole seconds
$ws.Range("V2").Formula = "=HOUR(U2)*3600+MINUTE(U2)*60+SECOND(U2)"
$ws.Range("V3:V13").Formula = "=HOUR(U3)*3600+MINUTE(U3)*60+SECOND(U3)"

# New "P" column: pace derived from the Second column
$ws.Range("W2").Formula = "=V2/42.195"
$ws.Range("W3:W13").Formula = "=V3/42.195"

# Apply number formatting (matches existing style used for numeric helper columns)
$ws.Range("V2:W13").NumberFormat = "0.00"

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection as left by the author
$ws.Range("W20").Select()
